$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

$ws.Range("C3").Value = "Todos los campos son llenados de forma correcta"
$ws.Range("D3").Value = "El examen se almacena en el sistema "
$ws.Range("B4").Value = "Numero preguntas invalido"
$ws.Range("B3").Value = "Examen creado con éxito"
$ws.Range("D4").Value = "Error: Numero de preguntas invalido"
$ws.Range("C4").Value = "Se ingresa un numero invalido para el total de preguntas"
$ws.Range("B5").Value = "Caracteres invalidos preguntas"
$ws.Range("C5").Value = "Se ingresan caracteres especiales en el total de preguntas"
$ws.Range("D5").Value = "Error: Ingresar un numero"
$ws.Range("B6").Value = "Letras en numero de preguntas"
$ws.Range("C6").Value = "Se escribren letas en el total de preguntas en vez de numeros"
$ws.Range("D6").Value = "Error:Ingresar un numero"

$ws.Range("B7").Select()
